$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (Title, Author, Date, Synopsis, URL)
$ws.Range("A1").Value = "Title"
$ws.Range("B1").Value = "Author"
$ws.Range("C1").Value = "Date"
$ws.Range("D1").Value = "Synopsis"
$ws.Range("E1").Value = "URL"

# Header formatting: bold, centered/top-aligned, thin box border
$a1 = $ws.Range("A1")
$a1.Borders.LineStyle = 1
$a1.Font.Bold = $true
$a1.HorizontalAlignment = -4108
$a1.VerticalAlignment = -4160
$a1.Copy()
$ws.Range("B1:E1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Article rows scraped so far (Author/Date/Synopsis not yet populated)
$ws.Range("A2").Value = "Patience pays off in quantum computing"
$ws.Range("E2").Value = "https://www.ft.com/content/63cf560a-69f3-42c2-82c3-7e8d97282716"
$ws.Range("A3").Value = "Quantum computing is struggling to reach its silicon moment"
$ws.Range("E3").Value = "https://www.ft.com/content/bfe5fff4-3d78-4bea-9e31-da8ca1f77151"
$ws.Range("A4").Value = "Quantum computing is overshadowed by rapid advances in AI"
$ws.Range("E4").Value = "https://www.ft.com/content/e3e2b721-9971-47b1-aa86-f210804ebc3e"
$ws.Range("A5").Value = "Letter: Quantum computing does the hard stuff AI alone can’t crack"
$ws.Range("E5").Value = "https://www.ft.com/content/bc1298fd-868c-4c21-b9c1-81716b583c7a"
$ws.Range("A6").Value = "Scientific breakthrough gives new hope to building quantum computers"
$ws.Range("E6").Value = "https://www.ft.com/content/f1d26918-67c5-4b11-b47b-4904606a002f"
$ws.Range("A7").Value = "Microsoft claims quantum breakthrough after 20-year pursuit of elusive particle"
$ws.Range("E7").Value = "https://www.ft.com/content/a60f44f5-81ca-4e66-8193-64c956b09820"
$ws.Range("A8").Value = "The mysterious promise of the quantum future"
$ws.Range("E8").Value = "https://www.ft.com/content/7ce2f04a-f8e2-4e3d-8602-776647c520c9"
$ws.Range("A9").Value = "Quantum computing breakthroughs draw investment back to sector"
$ws.Range("E9").Value = "https://www.ft.com/content/d0b486ab-ed6c-46f0-b7b6-66cc60780efe"
